# MOSIP-17570 added new supervisor rejected templates for SMS and EMAIL
#
# Adds 18 new rows (1716-1733) to Sheet1 describing two new template types
# (RPR_SUP_REJECT_EMAIL / RPR_SUP_REJECT_SMS) for each language, plus a
# RPR_SUP_REJECT_EMAIL_SUBJECT entry for each language.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# lang_code -> value used in column A
$eng = "eng"
$fra = "fra"
$ara = "ara"
$hin = "hin"
$kan = "kan"
$tam = "tam"

$emailCode = "RPR_SUP_REJECT_EMAIL"
$emailDesc = "Template for Supervisor Reject Email"
$smsCode = "RPR_SUP_REJECT_SMS"
$smsDesc = "Template for Supervisor Reject SMS"
$subjectCode = "RPR_SUP_REJECT_EMAIL_SUBJECT"
$subjectDesc = "Template for Supervisor Reject Email Subject"

# (row, lang, code, description)
$rows = @(
    @(1716, $eng, $emailCode, $emailDesc),
    @(1717, $eng, $smsCode,   $smsDesc),
    @(1718, $fra, $emailCode, $emailDesc),
    @(1719, $fra, $smsCode,   $smsDesc),
    @(1720, $ara, $emailCode, $emailDesc),
    @(1721, $ara, $smsCode,   $smsDesc),
    @(1722, $hin, $emailCode, $emailDesc),
    @(1723, $hin, $smsCode,   $smsDesc),
    @(1724, $kan, $emailCode, $emailDesc),
    @(1725, $kan, $smsCode,   $smsDesc),
    @(1726, $tam, $emailCode, $emailDesc),
    @(1727, $tam, $smsCode,   $smsDesc),
    @(1728, $eng, $subjectCode, $subjectDesc),
    @(1729, $fra, $subjectCode, $subjectDesc),
    @(1730, $ara, $subjectCode, $subjectDesc),
    @(1731, $hin, $subjectCode, $subjectDesc),
    @(1732, $kan, $subjectCode, $subjectDesc),
    @(1733, $tam, $subjectCode, $subjectDesc)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]

    # Column D always holds the literal text "TRUE" (shared with the rest of
    # the sheet). Assigning the string "TRUE"/"FALSE" directly makes Excel's
    # COM layer coerce it to a boolean cell, so instead copy the existing
    # literal-text "TRUE" cell from D2 and paste it into the new cell - this
    # preserves both the text cell type and the "D" column style.
    $ws.Range("D2").Copy() | Out-Null
    $ws.Range("D" + $rowNum).PasteSpecial() | Out-Null
}

$excel.CutCopyMode = $false

# Restore the view roughly to where the author left it: scrolled down near
# the bottom of the sheet with H1718 selected.
$excel.ActiveWindow.ScrollRow = 1711
$ws.Range("H1718").Select() | Out-Null
